$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "food super group"
$ws.Range("B6").Value = "condiment"

$ws.Range("A7").Select()
